# Optuna Attempt (go back with original)
# Applies updated forecast metrics to the "Forecast Comparison" sheet
# and refreshed rollup totals to the "Summary" sheet.

$wb = $excel.ActiveWorkbook

# --- Sheet 1: Forecast Comparison ---
$ws1 = $wb.Worksheets.Item("Forecast Comparison")

# Row 2
$ws1.Range("L2").Value = 1.18

# Row 3
$ws1.Range("D3").Value = 131
$ws1.Range("H3").Value = 3.44
$ws1.Range("L3").Value = 0.98

# Row 4
$ws1.Range("H4").Value = 3.47
$ws1.Range("L4").Value = 1.13

# Row 5
$ws1.Range("H5").Value = 2.43
$ws1.Range("L5").Value = 1.17

# Row 6
$ws1.Range("D6").Value = 96
$ws1.Range("H6").Value = 1.4
$ws1.Range("L6").Value = 1.17

# Row 7
$ws1.Range("D7").Value = 94
$ws1.Range("H7").Value = 0.4
$ws1.Range("I7").Value = "High"
$ws1.Range("J7").Value = "Urgent"
$ws1.Range("L7").Value = 1.11

# Row 8
$ws1.Range("D8").Value = 87
$ws1.Range("H8").Value = 0
$ws1.Range("I8").Value = "High"
$ws1.Range("L8").Value = 1.19

# Row 9
$ws1.Range("D9").Value = 89
$ws1.Range("L9").Value = 1.16

# Row 10
$ws1.Range("D10").Value = 91

# Row 11
$ws1.Range("D11").Value = 83
$ws1.Range("L11").Value = 1.18

# Row 12
$ws1.Range("D12").Value = 81
$ws1.Range("L12").Value = 1.04

# Row 13
$ws1.Range("L13").Value = 1.1

# Row 14
$ws1.Range("L14").Value = 1.08

# Row 15
$ws1.Range("L15").Value = 1.02

# Row 16
$ws1.Range("L16").Value = 0.88

# Row 17
$ws1.Range("L17").Value = 0.87

# --- Sheet 2: Summary ---
# Values in this column are stored as text (e.g. "1234"), not numbers, so
# force a text number format before assigning to avoid Excel auto-converting
# the numeric-looking strings into real numbers.
$ws2 = $wb.Worksheets.Item("Summary")

$ws2.Range("B9:B12").NumberFormat = "@"
$ws2.Range("B14").NumberFormat = "@"

$ws2.Range("B9").Value = "1456"
$ws2.Range("B10").Value = "809"
$ws2.Range("B11").Value = "441"
$ws2.Range("B12").Value = "131"
$ws2.Range("B14").Value = "72"
